$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6182220821905557
$ws.Range("B3").Value = 0.8222533184040642
$ws.Range("B4").Value = 0.6894054741978975
$ws.Range("B5").Value = 0.6894054741978975
$ws.Range("B6").Value = 0.6894054762660493
$ws.Range("B7").Value = 0.6345773421346755
